$wb = $excel.ActiveWorkbook

# --- Sheet "measurement" -> rename to "observation" and drop its Variance column ---
$obs = $wb.Worksheets.Item("measurement")
$obs.Name = "observation"

# The Variance column (B) moves out to its own sheet, so delete it here.
[void]$obs.Columns.Item(2).Delete()
[void]$obs.Range("B23").Select()

# --- New sheet "variance" holding the data that used to live in column B ---
$var = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $obs)
$var.Name = "variance"

$varianceValues = @(3, 4, 2, 5, 1, 8, 9, 1, 2, 3, 8, 4, 1)

$var.Range("A1").Value = "Variance"
for ($i = 0; $i -lt $varianceValues.Length; $i++) {
    $var.Cells.Item($i + 2, 1).Value = $varianceValues[$i]
}

# Make "variance" the active/selected sheet and cell, like in the saved file.
[void]$var.Range("E13").Select()
[void]$var.Select()
